$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Create the table (ListObject) over the full data range first.
# The header row (A1:U1) already carries explicit bold/fill/border direct
# formatting, so build the table on a throwaway range elsewhere and then
# resize it onto the real range - this avoids the engine capturing the
# pre-existing header formatting as an extra header-row dxf style.
$ws.Range("W70:Y72").Value = "tmp"
$tmpRange = $ws.Range("W70:Y72")
$listObj = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $tmpRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$listObj.Name = "Table1"
$listObj.Resize($ws.Range("A1:U64"))
$ws.Range("W70:Y72").Clear()
$listObj.TableStyle = ""

# Rename the "_old" / "_new" header suffixes to "_FV2404" / "_FV2410".
# Writing straight into the header cells also renames the matching
# ListColumns, since they're a live view over the header row.
$oldNamesFV2404 = @(
    "Segmentname_old","Segmentgruppe_old","Segment_old","Datenelement_old","Segment ID_old",
    "Code_old","Qualifier_old","Beschreibung_old","Bedingungsausdruck_old","Bedingung_old"
)
$newNamesFV2404 = @(
    "Segmentname_FV2404","Segmentgruppe_FV2404","Segment_FV2404","Datenelement_FV2404","Segment ID_FV2404",
    "Code_FV2404","Qualifier_FV2404","Beschreibung_FV2404","Bedingungsausdruck_FV2404","Bedingung_FV2404"
)
for ($i = 0; $i -lt $oldNamesFV2404.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $newNamesFV2404[$i]
}

# Column K ("diff") keeps its text, but still needs to be re-written so the
# table's 11th ListColumn name (currently the auto-generated placeholder
# from table creation) is refreshed to match the header cell.
$ws.Cells.Item(1, 11).Value = "diff"

$oldNamesFV2410 = @(
    "Segmentname_new","Segmentgruppe_new","Segment_new","Datenelement_new","Segment ID_new",
    "Code_new","Qualifier_new","Beschreibung_new","Bedingungsausdruck_new","Bedingung_new"
)
$newNamesFV2410 = @(
    "Segmentname_FV2410","Segmentgruppe_FV2410","Segment_FV2410","Datenelement_FV2410","Segment ID_FV2410",
    "Code_FV2410","Qualifier_FV2410","Beschreibung_FV2410","Bedingungsausdruck_FV2410","Bedingung_FV2410"
)
for ($i = 0; $i -lt $oldNamesFV2410.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $newNamesFV2410[$i]
}

# Freeze the header row (pane split after row 1)
$ws.Activate()
$ws.Application.ActiveWindow.SplitRow = 1
$ws.Application.ActiveWindow.FreezePanes = $true
